$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Summary")

# Insert 7 new columns starting at H (old H..N), shifting old H:O to O:V.
$ws.Columns("H:N").Insert()

# Populate the six new "Tier" header cells (H1:M1); N1 stays blank,
# matching the old spacer pattern that used to sit at G1.
$ws.Range("H1").Value = "Tier 1 `nCount"
$ws.Range("I1").Value = "Tier 2 `nCount"
$ws.Range("J1").Value = "Tier 3 `nCount"
$ws.Range("K1").Value = "Tier 4 `nCount"
$ws.Range("L1").Value = "Tier 5 `nCount"
$ws.Range("M1").Value = "Tier`n(unassigned) `nCount"

# Re-establish the AutoFilter over the full, now-wider header row.
$ws.AutoFilterMode = $false
$null = $ws.Range("A1:V1").AutoFilter()

# Update the hidden _FilterDatabase defined name to match the new extent.
foreach ($n in $wb.Names) {
    if ($n.Name -eq "Product Summary!_FilterDatabase") {
        $n.RefersTo = "='Product Summary'!`$A`$1:`$V`$1"
    }
}
